$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7, column A: new error message
$ws.Range("A7").Value = "Cannot find class [org.springframework.jdbc.datasource.DriverManagerDataSource] for bean with name 'dataSource'"

# Row 7, column B: the solution text (multi-line), word-wrapped
$b7text = @"
add <dependency>
    <groupId>org.springframework</groupId>
    <artifactId>spring-jdbc</artifactId>
    <version>3.0.3.RELEASE</version>
</dependency>  in pom
"@
$ws.Range("B7").Value = $b7text
$ws.Range("B7").WrapText = $true

# Taller row to fit the wrapped text
$ws.Rows.Item(7).RowHeight = 75

# Move/select the next cell, matching the saved workbook's cursor position
$null = $ws.Range("D7").Select()
